# Append five new subject rows (s16..s20) to the Simulation Manifest sheet,
# matching the "meltpatch" feedback records already present in rows 2-16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New records to append starting at row 17 (sheet currently has data in A1:I16).
$newRows = @(
    @("s16", "s16_IMG_3178.jpeg", "meltpatch", "1734", "1350", "104", "52", "168", "2"),
    @("s17", "s17_IMG_3177.jpeg", "meltpatch", "1104", "731",  "104", "52", "73",  "2"),
    @("s18", "s18_IMG_3174.jpeg", "meltpatch", "1679", "1816", "104", "52", "37",  "2"),
    @("s19", "s19_IMG_3179.jpeg", "meltpatch", "2615", "2426", "104", "52", "139", "2"),
    @("s20", "s20_IMG_3176.jpeg", "meltpatch", "902",  "1081", "104", "52", "113", "2")
)

$startRow = 17
$endRow = $startRow + $newRows.Length - 1

# Columns D:I hold numeric-looking values (x, y, toleranceA, toleranceB, theta,
# ratio) that must be kept as TEXT, matching every existing row in the sheet.
# Pre-format the target block as Text so assigning e.g. "1734" doesn't get
# silently reinterpreted as a number.
$ws.Range("D$startRow`:I$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $rowValues = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $rowValues[0]  # !subject_id
    $ws.Cells.Item($r, 2).Value = $rowValues[1]  # #file_name
    $ws.Cells.Item($r, 3).Value = $rowValues[2]  # #feedback_1_id
    $ws.Cells.Item($r, 4).Value = $rowValues[3]  # #feedback_1_x
    $ws.Cells.Item($r, 5).Value = $rowValues[4]  # #feedback_1_y
    $ws.Cells.Item($r, 6).Value = $rowValues[5]  # #feedback_1_toleranceA
    $ws.Cells.Item($r, 7).Value = $rowValues[6]  # #feedback_1_toleranceB
    $ws.Cells.Item($r, 8).Value = $rowValues[7]  # #feedback_1_theta
    $ws.Cells.Item($r, 9).Value = $rowValues[8]  # #minor_to_major_ratio
}
